# refatorando o consolidador para modelo ETL
# Update absenteeism data rows 2-11 with refreshed ETL-sourced values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 66229; B = "Lívia Rodrigues";         C = "TI";                      D = "Outros";              E = 5; F = 45106; G = 11700.68 }
    @{ Row = 3;  A = 67065; B = "Eloah Cardoso";            C = "Jurídico";                D = "Viagem de negócios";  E = 3; F = 45089; G = 7055.04 }
    @{ Row = 4;  A = 65431; B = "Sr. João Pedro Cardoso";   C = "Jurídico";                D = "Consulta médica";     E = 6; F = 45088; G = 12091.66 }
    @{ Row = 5;  A = 56930; B = "Davi Lucas Cardoso";       C = "Marketing";               D = "Problemas pessoais";  E = 4; F = 45092; G = 4451.09 }
    @{ Row = 6;  A = 58621; B = "Júlia da Conceição";       C = "Marketing";               D = "Problemas pessoais";  E = 7; F = 45082; G = 8222.4 }
    @{ Row = 7;  A = 35199; B = "Nicolas Viana";            C = "Financeiro";              D = "Doença";              E = 1; F = 45102; G = 11688.95 }
    @{ Row = 8;  A = 56573; B = "Ian Moreira";              C = "Financeiro";              D = "Problemas pessoais";  E = 5; F = 45092; G = 4904.81 }
    @{ Row = 9;  A = 25782; B = "Pietra Lopes";             C = "P&D";                     D = "Viagem de negócios";  E = 5; F = 45082; G = 3191.81 }
    @{ Row = 10; A = 47683; B = "Rodrigo Silveira";         C = "P&D";                     D = "Viagem de negócios";  E = 1; F = 45078; G = 10926.93 }
    @{ Row = 11; A = 43097; B = "Maria Monteiro";           C = "Atendimento ao Cliente";  D = "Problemas pessoais";  E = 7; F = 45104; G = 7432.14 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
}
